# Refresh the crypto price/volume snapshot pulled in by the GitHub Actions job.
# Numeric-looking price strings are entered with a leading apostrophe so Excel
# keeps them as text (matching the sheet's existing inline-string cells)
# instead of auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.958.90'
$ws.Range("E2").Value = '  -4.37%  '
$ws.Range("D3").Value = '2.606.10'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''516.26'
$ws.Range("E5").Value = '  -1.91%  '
$ws.Range("D6").Value = '''142.08'
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("D9").Value = '''6.69'
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("E10").Value = '  -3.04%  '
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").Value = '3.068.25'
$ws.Range("E13").Value = '  -3.36%  '
$ws.Range("D14").Value = '57.969.70'
$ws.Range("E14").Value = '  -4.31%  '
$ws.Range("D15").Value = '''20.81'
$ws.Range("E15").Value = '  -2.15%  '
$ws.Range("E16").Value = '  -1.94%  '
$ws.Range("D17").Value = '2.603.59'
$ws.Range("E17").Value = '  -4.17%  '
$ws.Range("E18").Value = '  -2.70%  '
$ws.Range("D19").Value = '''333.71'
$ws.Range("E19").Value = '  -3.53%  '
$ws.Range("D20").Value = '''10.31'
$ws.Range("E20").Value = '  -2.76%  '
$ws.Range("E21").Value = '  -3.29%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("E24").Value = '  -1.81%  '
$ws.Range("E25").Value = '  -2.50%  '
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("D27").Value = '''7.07'
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("E28").Value = '  -4.52%  '
$ws.Range("D29").Value = '''6.59'
$ws.Range("E29").Value = '  -3.49%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("E31").Value = '  -1.52%  '
$ws.Range("D32").Value = '''150.06'
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("E33").Value = '  -2.02%  '
$ws.Range("E34").Value = '  -4.08%  '
$ws.Range("E35").Value = '  -5.95%  '
$ws.Range("D36").Value = '''0.896'
$ws.Range("E36").Value = '  -4.51%  '
$ws.Range("D37").Value = '''36.60'
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("D38").Value = '''0.837'
$ws.Range("E38").Value = '  -4.31%  '
$ws.Range("E39").Value = '  -6.03%  '
$ws.Range("E40").Value = '  -1.95%  '
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("D43").Value = '''0.0962'
$ws.Range("E43").Value = '  -2.52%  '
$ws.Range("D44").Value = '''267.43'
$ws.Range("E44").Value = '  -5.51%  '
$ws.Range("D45").Value = '''10.61'
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("E46").Value = '  -4.90%  '
$ws.Range("D47").Value = '''0.0531'
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("D48").Value = '2.026.48'
$ws.Range("E48").Value = '  -5.58%  '
$ws.Range("E49").Value = '  -2.18%  '
$ws.Range("E50").Value = '  -4.71%  '
$ws.Range("D51").Value = '''18.16'
$ws.Range("E51").Value = '  -4.88%  '
